# acpt: add scenarios for Font.language_id
#
# Adds a 5th slide (blank layout) with three textboxes that exercise
# TextRange.LanguageID / Font.language_id:
#   Shape 0 - no explicit language set at all
#   Shape 1 - French  (fr-FR)
#   Shape 2 - Polish  (pl-PL)

$p = $ppt.ActivePresentation

# PowerPoint stores shape geometry in EMU (English Metric Units) in the
# OOXML, but the COM Left/Top/Width/Height (and AddTextbox) properties
# are expressed in points, where 1 pt = 12700 EMU.
$emuPerPt = 12700

$s = $p.Slides.Add($p.Slides.Count + 1, 12)  # ppLayoutBlank

# --- Shape 0: "Shape 0 - no explicit language id" -----------------------
$shp0 = $s.Shapes.AddTextbox(1, 2952093 / $emuPerPt, 1437501 / $emuPerPt, 3239814 / $emuPerPt, 369332 / $emuPerPt)
$shp0.TextFrame.WordWrap = $false
$shp0.TextFrame.AutoSize = 1
$shp0.Fill.Visible = 0
[void]$shp0.TextFrame.TextRange.InsertAfter("Shape 0 " + [char]0x2013 + " no explicit language id")

# --- Shape 1: "Shape 1 - MSO_LANGUAGE_ID.FRENCH" -------------------------
$shp1 = $s.Shapes.AddTextbox(1, 2613453 / $emuPerPt, 3244334 / $emuPerPt, 3917095 / $emuPerPt, 369332 / $emuPerPt)
$shp1.TextFrame.WordWrap = $false
$shp1.TextFrame.AutoSize = 1
$shp1.Fill.Visible = 0
[void]$shp1.TextFrame.TextRange.InsertAfter("Shape 1 " + [char]0x2013 + " MSO_LANGUAGE_ID.FRENCH")
$shp1.TextFrame.TextRange.LanguageID = "fr-FR"

# --- Shape 2: "Shape 2 - MSO_LANGUAGE_ID.POLISH" -------------------------
$shp2 = $s.Shapes.AddTextbox(1, 2654874 / $emuPerPt, 5051167 / $emuPerPt, 3834253 / $emuPerPt, 369332 / $emuPerPt)
$shp2.TextFrame.WordWrap = $false
$shp2.TextFrame.AutoSize = 1
$shp2.Fill.Visible = 0
[void]$shp2.TextFrame.TextRange.InsertAfter("Shape 2 " + [char]0x2013 + " MSO_LANGUAGE_ID.POLISH")
$shp2.TextFrame.TextRange.LanguageID = "pl-PL"
